$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: insert two blank columns before column D ---
# (existing columns D:K, with 8 quarters of data, shift right to become F:M)
$ws.Columns("D:E").Insert()

# --- Step 2: copy number formatting from the (now-shifted) original D column
#     (currently sitting in column F) into the two new D:E columns, one contiguous
#     block of rows at a time, so the new cells inherit the right date/number style ---
$ws.Range("F7:F35").Copy()
$ws.Range("D7:D35").PasteSpecial(-4122)
$ws.Range("F7:F35").Copy()
$ws.Range("E7:E35").PasteSpecial(-4122)
$ws.Range("F38:F77").Copy()
$ws.Range("D38:D77").PasteSpecial(-4122)
$ws.Range("F38:F77").Copy()
$ws.Range("E38:E77").PasteSpecial(-4122)
$ws.Range("F80:F102").Copy()
$ws.Range("D80:D102").PasteSpecial(-4122)
$ws.Range("F80:F102").Copy()
$ws.Range("E80:E102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Step 3: write the two new quarters of data (period ending 2019-01-31 and
#     2018-10-29) into the new D and E columns ---
$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("D8").Value = 228600
$ws.Range("E8").Value = 225400
$ws.Range("D9").Value = 7400
$ws.Range("E9").Value = 7400
$ws.Range("D10").Value = 221200
$ws.Range("E10").Value = 218000
$ws.Range("D12").Value = 57300
$ws.Range("E12").Value = 44700
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = "NA"
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 0
$ws.Range("D17").Value = 117000
$ws.Range("E17").Value = 100200
$ws.Range("D18").Value = 111600
$ws.Range("E18").Value = 125200
$ws.Range("D20").Value = 4800
$ws.Range("E20").Value = 3800
$ws.Range("D21").Value = 121300
$ws.Range("E21").Value = 130700
$ws.Range("D22").Value = "NA"
$ws.Range("E22").Value = 0
$ws.Range("D23").Value = 116400
$ws.Range("E23").Value = 129000
$ws.Range("D24").Value = -243700
$ws.Range("E24").Value = 2300
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = 360100
$ws.Range("E26").Value = 126600
$ws.Range("D27").Value = 360100
$ws.Range("E27").Value = 126600
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = 0
$ws.Range("E29").Value = 0
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = -4800
$ws.Range("E32").Value = -3800
$ws.Range("D33").Value = 360100
$ws.Range("E33").Value = 126600
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = 360100
$ws.Range("E35").Value = 126600
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("D41").Value = 314800
$ws.Range("E41").Value = 353600
$ws.Range("D42").Value = 378600
$ws.Range("E42").Value = 281000
$ws.Range("D43").Value = 178800
$ws.Range("E43").Value = 129100
$ws.Range("D44").Value = 9800
$ws.Range("E44").Value = 10400
$ws.Range("D45").Value = 15000
$ws.Range("E45").Value = 12800
$ws.Range("D46").Value = 897000
$ws.Range("E46").Value = 787000
$ws.Range("D47").Value = 157200
$ws.Range("E47").Value = 114100
$ws.Range("D48").Value = 56800
$ws.Range("E48").Value = 57000
$ws.Range("D49").Value = 63700
$ws.Range("E49").Value = 63700
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 247600
$ws.Range("E52").Value = 2600
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 1422300
$ws.Range("E54").Value = 1024400
$ws.Range("D57").Value = 10900
$ws.Range("E57").Value = 8700
$ws.Range("D58").Value = 0
$ws.Range("E58").Value = 0
$ws.Range("D59").Value = 94500
$ws.Range("E59").Value = 82200
$ws.Range("D60").Value = 105500
$ws.Range("E60").Value = 90900
$ws.Range("D61").Value = 100
$ws.Range("E61").Value = 100
$ws.Range("D62").Value = 29300
$ws.Range("E62").Value = 17400
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 134800
$ws.Range("E66").Value = 108400
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = -880400
$ws.Range("E72").Value = -1240500
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 1287500
$ws.Range("E76").Value = 916000
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("D81").Value = 360100
$ws.Range("E81").Value = 126600
$ws.Range("D83").Value = 4900
$ws.Range("E83").Value = 1700
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = 104600
$ws.Range("E89").Value = 164600
$ws.Range("D91").Value = -2900
$ws.Range("E91").Value = -15200
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = -142800
$ws.Range("E94").Value = -63300
$ws.Range("D96").Value = 0
$ws.Range("E96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = -1100
$ws.Range("E100").Value = 3900
$ws.Range("D101").Value = 0
$ws.Range("E101").Value = 0
$ws.Range("D102").Value = -39400
$ws.Range("E102").Value = 105200

# --- Step 4: a handful of historical figures were also corrected as part of this
#     refresh; apply those on top of the shifted columns ---
$ws.Range("G8").Value = 213700
$ws.Range("G10").Value = 208100
$ws.Range("G17").Value = 97400
$ws.Range("H61").Value = 14500
$ws.Range("H62").Value = 240600
$ws.Range("I91").Value = -1200
$ws.Range("J91").Value = -1500
$ws.Range("H94").Value = -18300
$ws.Range("I94").Value = -26600
$ws.Range("H102").Value = 34300
$ws.Range("I102").Value = 14100
